$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: add empty cells X15:AM15 with style matching s="3" (copy format from A67)
$ws.Range("A67").Copy() | Out-Null
$ws.Range("X15:AM15").PasteSpecial(-4122) | Out-Null

# Rows 16-30: add empty cell X{row} with style matching s="3"
$ws.Range("A67").Copy() | Out-Null
$ws.Range("X16:X30").PasteSpecial(-4122) | Out-Null

# Rows 68-75: new data rows with values
$ws.Cells.Item(68, 1).Value = 0.2348168041
$ws.Cells.Item(68, 2).Value = 0.1932066166
$ws.Cells.Item(68, 3).Value = 0.2024011192
$ws.Cells.Item(68, 4).Value = 1.0
$ws.Cells.Item(68, 5).Value = 0.4275661366
$ws.Cells.Item(68, 6).Value = 0.1545882395
$ws.Cells.Item(68, 7).Value = 0.3663011142
$ws.Cells.Item(68, 8).Value = 0.1355981002
$ws.Cells.Item(68, 9).Value = 0.1716433163
$ws.Cells.Item(68, 10).Value = 0.05840731859
$ws.Cells.Item(68, 11).Value = 0.1815098907
$ws.Cells.Item(68, 12).Value = 0.2148618783
$ws.Cells.Item(68, 13).Value = 0.0804714633
$ws.Cells.Item(68, 14).Value = 0.07956648866
$ws.Cells.Item(68, 15).Value = 0.08492013035
$ws.Cells.Item(68, 16).Value = 0.1206863846
$ws.Cells.Item(68, 17).Value = 0.106633975
$ws.Cells.Item(68, 18).Value = 0.07196264847
$ws.Cells.Item(68, 19).Value = 0.07154999437
$ws.Cells.Item(68, 20).Value = 0.7228005764
$ws.Cells.Item(68, 21).Value = 0.3306782996
$ws.Cells.Item(68, 22).Value = 0.0

$ws.Cells.Item(69, 1).Value = 0.03622695895
$ws.Cells.Item(69, 2).Value = 0.05903942799
$ws.Cells.Item(69, 3).Value = 0.3457125212
$ws.Cells.Item(69, 4).Value = 0.5572050479
$ws.Cells.Item(69, 5).Value = 0.424166329
$ws.Cells.Item(69, 6).Value = 0.1906959285
$ws.Cells.Item(69, 7).Value = 1.0
$ws.Cells.Item(69, 8).Value = 0.1900041546
$ws.Cells.Item(69, 9).Value = 0.04994227871
$ws.Cells.Item(69, 10).Value = 0.1670153272
$ws.Cells.Item(69, 11).Value = 0.3610296567
$ws.Cells.Item(69, 12).Value = 0.072531699
$ws.Cells.Item(69, 13).Value = 0.04344478551
$ws.Cells.Item(69, 14).Value = 0.04090720862
$ws.Cells.Item(69, 15).Value = 0.02951475714
$ws.Cells.Item(69, 16).Value = 0.1385135777
$ws.Cells.Item(69, 17).Value = 0.1245111833
$ws.Cells.Item(69, 18).Value = 0.1718209014
$ws.Cells.Item(69, 19).Value = 0.07125407247
$ws.Cells.Item(69, 20).Value = 0.3915861378
$ws.Cells.Item(69, 21).Value = 0.07874547971
$ws.Cells.Item(69, 22).Value = 0.0

$ws.Cells.Item(70, 1).Value = 0.0488832523
$ws.Cells.Item(70, 2).Value = 0.03363707051
$ws.Cells.Item(70, 3).Value = 0.1751917358
$ws.Cells.Item(70, 4).Value = 0.3278934577
$ws.Cells.Item(70, 5).Value = 0.1393134209
$ws.Cells.Item(70, 6).Value = 0.09500132692
$ws.Cells.Item(70, 7).Value = 1.0
$ws.Cells.Item(70, 8).Value = 0.2800308144
$ws.Cells.Item(70, 9).Value = 0.2139728105
$ws.Cells.Item(70, 10).Value = 0.1902788298
$ws.Cells.Item(70, 11).Value = 0.5086038764
$ws.Cells.Item(70, 12).Value = 0.3507530103
$ws.Cells.Item(70, 13).Value = 0.1986300702
$ws.Cells.Item(70, 14).Value = 0.01873973458
$ws.Cells.Item(70, 15).Value = 0.009034634675
$ws.Cells.Item(70, 16).Value = 0.05240802419
$ws.Cells.Item(70, 17).Value = 0.04692523531
$ws.Cells.Item(70, 18).Value = 0.1060700534
$ws.Cells.Item(70, 19).Value = 0.06657526051
$ws.Cells.Item(70, 20).Value = 0.2279955318
$ws.Cells.Item(70, 21).Value = 0.1003939348
$ws.Cells.Item(70, 22).Value = 0.0

$ws.Cells.Item(71, 1).Value = 0.04843373506
$ws.Cells.Item(71, 2).Value = 0.1248203532
$ws.Cells.Item(71, 3).Value = 0.08879431562
$ws.Cells.Item(71, 4).Value = 1.0
$ws.Cells.Item(71, 5).Value = 0.1304157129
$ws.Cells.Item(71, 6).Value = 0.1032157881
$ws.Cells.Item(71, 7).Value = 0.9790417206
$ws.Cells.Item(71, 8).Value = 0.1473743801
$ws.Cells.Item(71, 9).Value = 0.06723048184
$ws.Cells.Item(71, 10).Value = 0.06737927566
$ws.Cells.Item(71, 11).Value = 0.4210878534
$ws.Cells.Item(71, 12).Value = 0.09718495159
$ws.Cells.Item(71, 13).Value = 0.09207549049
$ws.Cells.Item(71, 14).Value = 0.04390416475
$ws.Cells.Item(71, 15).Value = 0.03722126539
$ws.Cells.Item(71, 16).Value = 0.04702293593
$ws.Cells.Item(71, 17).Value = 0.05871873036
$ws.Cells.Item(71, 18).Value = 0.2965884774
$ws.Cells.Item(71, 19).Value = 0.172424024
$ws.Cells.Item(71, 20).Value = 0.1665995265
$ws.Cells.Item(71, 21).Value = 0.359643318
$ws.Cells.Item(71, 22).Value = 0.0

$ws.Cells.Item(72, 1).Value = 0.2051235585
$ws.Cells.Item(72, 2).Value = 0.3372615558
$ws.Cells.Item(72, 3).Value = 0.09814005263
$ws.Cells.Item(72, 4).Value = 0.5368386057
$ws.Cells.Item(72, 5).Value = 1.0
$ws.Cells.Item(72, 6).Value = 0.4124115094
$ws.Cells.Item(72, 7).Value = 0.4334421579
$ws.Cells.Item(72, 8).Value = 0.2236848109
$ws.Cells.Item(72, 9).Value = 0.2889681047
$ws.Cells.Item(72, 10).Value = 0.1428138365
$ws.Cells.Item(72, 11).Value = 0.2358964746
$ws.Cells.Item(72, 12).Value = 0.1683372654
$ws.Cells.Item(72, 13).Value = 0.1252625801
$ws.Cells.Item(72, 14).Value = 0.07416808378
$ws.Cells.Item(72, 15).Value = 0.0986841733
$ws.Cells.Item(72, 16).Value = 0.07924967837
$ws.Cells.Item(72, 17).Value = 0.05365527954
$ws.Cells.Item(72, 18).Value = 0.1467018315
$ws.Cells.Item(72, 19).Value = 0.1135073523
$ws.Cells.Item(72, 20).Value = 0.8160838129
$ws.Cells.Item(72, 21).Value = 0.1787016218
$ws.Cells.Item(72, 22).Value = 0.0

$ws.Cells.Item(73, 1).Value = 0.06818114228
$ws.Cells.Item(73, 2).Value = 0.1346605553
$ws.Cells.Item(73, 3).Value = 0.2678248967
$ws.Cells.Item(73, 4).Value = 0.5929944419
$ws.Cells.Item(73, 5).Value = 0.4917079222
$ws.Cells.Item(73, 6).Value = 0.5040074283
$ws.Cells.Item(73, 7).Value = 1.0
$ws.Cells.Item(73, 8).Value = 0.4366139557
$ws.Cells.Item(73, 9).Value = 0.0888001289
$ws.Cells.Item(73, 10).Value = 0.04939210819
$ws.Cells.Item(73, 11).Value = 0.2510522887
$ws.Cells.Item(73, 12).Value = 0.3331240756
$ws.Cells.Item(73, 13).Value = 0.08169242136
$ws.Cells.Item(73, 14).Value = 0.02786573576
$ws.Cells.Item(73, 15).Value = 0.2101597203
$ws.Cells.Item(73, 16).Value = 0.0970351349
$ws.Cells.Item(73, 17).Value = 0.07849309195
$ws.Cells.Item(73, 18).Value = 0.143348789
$ws.Cells.Item(73, 19).Value = 0.06384451826
$ws.Cells.Item(73, 20).Value = 0.3821366544
$ws.Cells.Item(73, 21).Value = 0.1328850544
$ws.Cells.Item(73, 22).Value = 0.0

$ws.Cells.Item(74, 1).Value = 0.002209853287
$ws.Cells.Item(74, 2).Value = 0.002963303466
$ws.Cells.Item(74, 3).Value = 0.004229586931
$ws.Cells.Item(74, 4).Value = 0.02111864854
$ws.Cells.Item(74, 5).Value = 0.01767931857
$ws.Cells.Item(74, 6).Value = 0.3836951109
$ws.Cells.Item(74, 7).Value = 0.7145434297
$ws.Cells.Item(74, 8).Value = 0.7277384967
$ws.Cells.Item(74, 9).Value = 0.008136134816
$ws.Cells.Item(74, 10).Value = 0.004301396587
$ws.Cells.Item(74, 11).Value = 1.0
$ws.Cells.Item(74, 12).Value = 0.1477001322
$ws.Cells.Item(74, 13).Value = 0.458839411
$ws.Cells.Item(74, 14).Value = 0.04921012868
$ws.Cells.Item(74, 15).Value = 0.1175007274
$ws.Cells.Item(74, 16).Value = 0.02232443358
$ws.Cells.Item(74, 17).Value = 0.001885427626
$ws.Cells.Item(74, 18).Value = 0.05976978915
$ws.Cells.Item(74, 19).Value = 0.02916356813
$ws.Cells.Item(74, 20).Value = 0.02767817369
$ws.Cells.Item(74, 21).Value = 0.002621222777
$ws.Cells.Item(74, 22).Value = 0.0

$ws.Cells.Item(75, 1).Value = 0.09278163127
$ws.Cells.Item(75, 2).Value = 0.1035817937
$ws.Cells.Item(75, 3).Value = 0.3961774144
$ws.Cells.Item(75, 4).Value = 0.7902471368
$ws.Cells.Item(75, 5).Value = 0.9249079131
$ws.Cells.Item(75, 6).Value = 0.3575379968
$ws.Cells.Item(75, 7).Value = 0.5946909029
$ws.Cells.Item(75, 8).Value = 0.1246131888
$ws.Cells.Item(75, 9).Value = 0.2036263852
$ws.Cells.Item(75, 10).Value = 0.1157371923
$ws.Cells.Item(75, 11).Value = 1.0
$ws.Cells.Item(75, 12).Value = 0.1411153595
$ws.Cells.Item(75, 13).Value = 0.3169310008
$ws.Cells.Item(75, 14).Value = 0.09374102588
$ws.Cells.Item(75, 15).Value = 0.1328719086
$ws.Cells.Item(75, 16).Value = 0.2337012287
$ws.Cells.Item(75, 17).Value = 0.2487162114
$ws.Cells.Item(75, 18).Value = 0.309378615
$ws.Cells.Item(75, 19).Value = 0.2271617479
$ws.Cells.Item(75, 20).Value = 0.2475636502
$ws.Cells.Item(75, 21).Value = 0.06422184813
$ws.Cells.Item(75, 22).Value = 0.0

# Apply styles for rows 68-75 per column groups
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A68:L68").PasteSpecial(-4122) | Out-Null
$ws.Range("A101").Copy() | Out-Null
$ws.Range("M68:U68").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("V68").PasteSpecial(-4122) | Out-Null

$ws.Range("A15").Copy() | Out-Null
$ws.Range("A69:V69").PasteSpecial(-4122) | Out-Null

$ws.Range("A15").Copy() | Out-Null
$ws.Range("A70:V70").PasteSpecial(-4122) | Out-Null

$ws.Range("A15").Copy() | Out-Null
$ws.Range("A71:N71").PasteSpecial(-4122) | Out-Null
$ws.Range("A101").Copy() | Out-Null
$ws.Range("O71:Q71").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("R71:V71").PasteSpecial(-4122) | Out-Null

$ws.Range("A15").Copy() | Out-Null
$ws.Range("A72:N72").PasteSpecial(-4122) | Out-Null
$ws.Range("A101").Copy() | Out-Null
$ws.Range("O72:Q72").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("R72:V72").PasteSpecial(-4122) | Out-Null

$ws.Range("A15").Copy() | Out-Null
$ws.Range("A73:N73").PasteSpecial(-4122) | Out-Null
$ws.Range("A101").Copy() | Out-Null
$ws.Range("O73:Q73").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("R73:V73").PasteSpecial(-4122) | Out-Null

$ws.Range("A15").Copy() | Out-Null
$ws.Range("A74:N74").PasteSpecial(-4122) | Out-Null
$ws.Range("A101").Copy() | Out-Null
$ws.Range("O74:Q74").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("R74:V74").PasteSpecial(-4122) | Out-Null

$ws.Range("A15").Copy() | Out-Null
$ws.Range("A75:N75").PasteSpecial(-4122) | Out-Null
$ws.Range("A101").Copy() | Out-Null
$ws.Range("O75:Q75").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("R75:V75").PasteSpecial(-4122) | Out-Null

# Row 91: add empty cells Y91:AG91 with style matching s="4"
$ws.Range("A101").Copy() | Out-Null
$ws.Range("Y91:AG91").PasteSpecial(-4122) | Out-Null

# Rows 94-98: add empty cells AA:AC with style matching s="4"
$ws.Range("A101").Copy() | Out-Null
$ws.Range("AA94:AC98").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0